# Refresh workbook with FDIC data and peer medians
# Updates shares outstanding (10-Q refresh) and converts several
# derived formula cells into their recalculated static values,
# plus refreshes the narrative assessment text on three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Assumptions
# ---------------------------------------------------------------
$wsAssump = $wb.Worksheets.Item("Assumptions")
$wsAssump.Range("B11").Value = 69.34

# ---------------------------------------------------------------
# Sheet: Base Case
# ---------------------------------------------------------------
$wsBase = $wb.Worksheets.Item("Base Case")
$wsBase.Range("C7").Value = 24.7
$wsBase.Range("C9").Value = 48.04
$wsBase.Range("C11").Value = 38.43
$wsBase.Range("C13").Value = 20.1
$wsBase.Range("C15").Value = 13.15
$wsBase.Range("C17").Value = 265

$wsBase.Range("A20").Value = 'Base case NCO normalization results in 20.1 bps CET1 burn.'
$wsBase.Range("A21").Value = 'New cushion of 265 bps ($38.43M) remains above the 10.5% buffer but far tighter than prior analysis.'
$wsBase.Range("A22").Value = 'Dividend ($77M annual) remains serviceable; share repurchases require tighter pacing.'
$wsBase.Range("A23").Value = 'Rating Impact: WATCHLIST - capital headroom compressed vs. prior 615 bps cushion.'

# ---------------------------------------------------------------
# Sheet: Industrial_Warehouse_Stress
# ---------------------------------------------------------------
$wsStress = $wb.Worksheets.Item("Industrial_Warehouse_Stress")

# Scenario 1: Base Case (5% cumulative loss rate)
$wsStress.Range("C15").Value = 96.55
$wsStress.Range("C17").Value = 77.24
$wsStress.Range("C19").Value = 40.4
$wsStress.Range("C21").Value = 12.95
$wsStress.Range("C22").Value = 245

$wsStress.Range("A24").Value = 'ASSESSMENT: Base case 5% loss = 40.4 bps CET1 burn. Cushion now 245 bps.'

# Scenario 2: Bear Case (15% cumulative loss rate)
$wsStress.Range("C31").Value = 289.65
$wsStress.Range("C33").Value = 231.72
$wsStress.Range("C35").Value = 121.2
$wsStress.Range("C37").Value = 12.14
$wsStress.Range("C38").Value = 164
$wsStress.Range("C39").Value = 69.34
$wsStress.Range("C40").Value = 3.34
$wsStress.Range("C42").Value = 32.82
$wsStress.Range("C43").Value = -9.2

$wsStress.Range("A45").Value = 'ASSESSMENT: Bear case 15% loss = 121.2 bps CET1 burn. Cushion compresses to 164 bps.'
$wsStress.Range("A46").Value = 'Combined with through-cycle NCO normalization (20.1 bps), total stress = 141.3 bps burn.'
$wsStress.Range("A47").Value = 'This would reduce CET1 to 12.14% (buffer = 164 bps) and TBVPS to $32.82 (-9.2% decline).'
